# refatoração - cálculos de apoio médio
# Adds std/min/max breakdown columns for "arrecadado", "apoio" and "contribuicoes",
# renames a couple of headers, and inserts the new statistics into the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) - final column layout A:V
# ---------------------------------------------------------------------------
$headers = @{
    "A1" = "modalidade"
    "B1" = "origem"
    "C1" = "total"
    "D1" = "total_sucesso"
    "E1" = "particip"
    "F1" = "taxa_sucesso"
    "G1" = "arrecadado_sucesso"
    "H1" = "arrecadado_avg"
    "I1" = "arrecadado_std"
    "J1" = "arrecadado_min"
    "K1" = "arrecadado_max"
    "L1" = "apoio_medio"
    "M1" = "apoio_std"
    "N1" = "apoio_min"
    "O1" = "apoio_max"
    "P1" = "contribuicoes"
    "Q1" = "contribuicoes_med"
    "R1" = "contribuicoes_std"
    "S1" = "contribuicoes_min"
    "T1" = "contribuicoes_max"
    "U1" = "menor_ano"
    "V1" = "maior_ano"
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# New header cells (Q1:V1) need the same bold/centered/bordered look as the
# rest of the header row - copy that formatting over from an existing header cell.
$ws.Range("A1:F1").Copy()
$ws.Range("Q1:V1").PasteSpecial(-4122) # xlPasteFormats

# ---------------------------------------------------------------------------
# Data rows (row 2 = sub/apoia.se, row 3 = sub/catarse)
# ---------------------------------------------------------------------------
$row2 = @{
    "A2" = "sub"
    "B2" = "apoia.se"
    "C2" = 627
    "D2" = 137
    "E2" = 0.9166666666666666
    "F2" = 0.2185007974481659
    "G2" = 39550.43984210649
    "H2" = 288.6893419131861
    "I2" = 682.4025885496077
    "J2" = 1.087396962410123
    "K2" = 5087.076865717208
    "L2" = 20.8884993069937
    "M2" = 15.28236810307352
    "N2" = 1.011042153300025
    "O2" = 84.0771316599004
    "P2" = 2063
    "Q2" = 15.05839416058394
    "R2" = 33.43095747941649
    "S2" = 1
    "T2" = 208
    "U2" = 2016
    "V2" = 2023
}

$row3 = @{
    "A3" = "sub"
    "B3" = "catarse"
    "C3" = 57
    "D3" = 15
    "E3" = 0.08333333333333333
    "F3" = 0.2631578947368421
    "G3" = 3636.517912678314
    "H3" = 242.4345275118876
    "I3" = 198.3989605548985
    "J3" = 10.98162164796783
    "K3" = 538.4389998789497
    "L3" = 24.89101283785623
    "M3" = 12.20279843561647
    "N3" = 10.98162164796783
    "O3" = 48.38432860277894
    "P3" = 145
    "Q3" = 9.666666666666666
    "R3" = 8.295150620062532
    "S3" = 1
    "T3" = 30
    "U3" = 2019
    "V3" = 2023
}

foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}
foreach ($addr in $row3.Keys) {
    $ws.Range($addr).Value = $row3[$addr]
}

# ---------------------------------------------------------------------------
# Number formats - match the existing column conventions:
#   "#,##0"       -> whole-number counts (style used by C:D, and now P:T)
#   "0.00%"       -> percentages (E:F, unchanged)
#   "R$ #,##0.00" -> currency (G:L, and now also M:O)
# ---------------------------------------------------------------------------
$ws.Range("G2:L3").NumberFormat = "R$ #,##0.00"
$ws.Range("M2:O3").NumberFormat = "R$ #,##0.00"
$ws.Range("C2:D3").NumberFormat = "#,##0"
$ws.Range("P2:T3").NumberFormat = "#,##0"
$ws.Range("E2:F3").NumberFormat = "0.00%"
